# ---------------------------------------------------------------------------
# Add 2022-Q3 data:
#  1. Insert a new sheet "2022-Q3" right before "2022-Q2" (i.e. as the 2nd
#     sheet) and populate it with the per-fund holdings detail.
#  2. Insert a new row at the top of the "总计" (summary) sheet's data with
#     the 2022-Q3 aggregate numbers, pushing the existing rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert new row for 2022-Q3
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift existing data rows (2..8) down to (3..9), bottom-up so we don't
# clobber a row before it has been read.
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Range("A$dest").Value = $summary.Range("A$r").Value2
    $summary.Range("B$dest").Value = $summary.Range("B$r").Value2
    $summary.Range("C$dest").Value = $summary.Range("C$r").Value2
    $summary.Range("D$dest").Value = $summary.Range("D$r").Value2
}

# New first data row: 2022-Q3
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 16
$summary.Range("D2").Value = 9.46

# Row 9 is brand-new (sheet grew from 8 to 9 data-bearing rows): give its "A"
# cell the same formatting (bold/border/centered) the other index cells use.
$summary.Range("A2").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Part 2: new "2022-Q3" worksheet with per-fund detail
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "008545", "泓德丰润三年持有期混合",             "74.71", "88.42", "4.05", "3.0258", 8),
    @(1,  "005395", "泓德臻远回报灵活配置混合",             "29.76", "93.35", "5.17", "1.5386", 7),
    @(2,  "001500", "泓德远见回报混合",                    "21.21", "92.75", "6.85", "1.4529", 7),
    @(3,  "010864", "泓德卓远混合A",                       "22.84", "92.87", "4.92", "1.1237", 5),
    @(4,  "004965", "泓德致远混合A",                       "16.19", "46.90", "4.48", "0.7253", 6),
    @(5,  "010865", "泓德卓远混合C",                       "10.32", "92.87", "4.92", "0.5077", 5),
    @(6,  "483003", "工银精选平衡混合",                    "15.78", "65.64", "2.48", "0.3913", 8),
    @(7,  "009447", "财通资管科技创新一年定期开放混合",      "8.69",  "94.05", "4.04", "0.3511", 8),
    @(8,  "004966", "泓德致远混合C",                       "2.54",  "46.90", "4.48", "0.1138", 6),
    @(9,  "516620", "国泰中证影视主题ETF",                 "0.94",  "99.07", "6.50", "0.0611", 5),
    @(10, "290012", "泰信行业精选灵活配置混合A",            "0.75",  "91.96", "7.17", "0.0538", 1),
    @(11, "159855", "银华中证影视主题ETF",                 "0.84",  "96.84", "6.33", "0.0532", 5),
    @(12, "001798", "泰康新回报灵活配置混合A",              "1.03",  "80.70", "3.59", "0.0370", 10),
    @(13, "517500", "国泰中证沪港深动漫游戏ETF",            "0.53",  "92.78", "3.31", "0.0175", 9),
    @(14, "001799", "泰康新回报灵活配置混合C",              "0.14",  "80.70", "3.59", "0.0050", 10),
    @(15, "002583", "泰信行业精选灵活配置混合C",            "0.04",  "91.96", "7.17", "0.0029", 1)
)

$r = 2
foreach ($row in $rows) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = "'" + $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = "'" + $row[3]
    $q3.Range("E$r").Value = "'" + $row[4]
    $q3.Range("F$r").Value = "'" + $row[5]
    $q3.Range("G$r").Value = "'" + $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Header styling (bold / bordered / centered) to match the other quarter
# sheets' header rows.
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Column A index cells use the same bold/bordered style too.
$q2.Range("A2").Copy()
$q3.Range("A2:A17").PasteSpecial(-4122)
